# ADD results from server
# Update row 2 (data row) values for each yearly sheet (2025, 2030, 2035, 2040, 2045, 2050)
# with new results received from the server.

$wb = $excel.ActiveWorkbook

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O")

$sheetsData = @(
    @{ Name = "2025"; Values = @(0, 1037.265132737054, 0, 0, 28926.05393052954, 0, 8095.925712661834, 0, 16171.06685703679, 0, 0, 48492.22142001599, 10595.37713982, 7083.519888241757, 6997.091525493762) },
    @{ Name = "2030"; Values = @(0, 4157.588990853394, 0, 0, 45991.90904307188, 0, 8095.925712661834, 0, 37079.12819938764, 0, 0, 54844.03303316472, 17449.04999683176, 9040.000118222546, 9731.320038188689) },
    @{ Name = "2035"; Values = @(2754.31755456332, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 13055.68800301272, 12867.93508263103) },
    @{ Name = "2040"; Values = @(2754.31755456332, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 13173.73845336726, 12867.93508263103) },
    @{ Name = "2045"; Values = @(5713.151062849596, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 13625.07956842681, 14949.88392520632) },
    @{ Name = "2050"; Values = @(5713.151062849596, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 14949.88392520632, 14949.88392520632) }
)

foreach ($entry in $sheetsData) {
    $sheetName = [string]$entry.Name
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $entry.Values
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])2"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
